# Generate Report for Handback
#
# Re-generates the handback status report: the
# "1bbacf06-3ffc-4215-b44d-39c6a6eca9d1" source file has dropped out of
# the current handback batch, so its row is removed from every sheet,
# and the handoff/handback timestamps for the still-present
# "075d7a73-6826-41de-8bfe-e6cf8e17f9ec" source file are refreshed on
# the locale sheets.

$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# "Overview" sheet: drop row 3 and rebuild the one surviving
# hyperlink (the engine's Range.Hyperlinks.Delete() clears every
# hyperlink on the sheet rather than scoping to the given range, so
# hyperlinks are rebuilt from known-good values after the wipe).
# -----------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Rows.Item(3).Delete()

$overview.Hyperlinks.Delete()
$overview.Hyperlinks.Add(
    $overview.Range("A2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md",
    "",
    "",
    "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"
)

# -----------------------------------------------------------------
# Locale sheets "zh-cn" / "de-de": drop row 3, refresh the
# "Correspond Handoff Datetime" (E2) / "Correspond Handback DateTime"
# (H2) cells, and rebuild the row-2 hyperlinks (A2, B2, D2, F2, G2).
# -----------------------------------------------------------------
$localeSheets = @(
    @{
        Name = "zh-cn"
        HandoffTime = "2016-03-19 10:38:09"
        HandbackTime = "2016-03-19 10:38:28"
        Links = @(
            @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md" },
            @{ Cell = "B2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"; Display = ".md" },
            @{ Cell = "D2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0da20fdcb2095ad1d9efce3b662a47508bef0b76/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf" },
            @{ Cell = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/f9d48f766536a39bd95aa2caf23c0530feaf8640/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md" },
            @{ Cell = "G2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/580a580d4fa62df7e11ef69dcbb94f2df72a474a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.zh-cn.xlf" }
        )
    },
    @{
        Name = "de-de"
        HandoffTime = "2016-03-19 10:38:12"
        HandbackTime = "2016-03-19 10:38:33"
        Links = @(
            @{ Cell = "A2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md" },
            @{ Cell = "B2"; Address = "https://github.com/OpenLocalizationTest/oltest/blob/7f983055ada23edc3bbdc6dcba8544a48464a646/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"; Display = ".md" },
            @{ Cell = "D2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/b037c97e912ee6d80f7a96fda6f39ce3f8dc28d4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf" },
            @{ Cell = "F2"; Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/3715e864f0ffc4408a68ce04248c2c8760d6980f/e2e/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.md" },
            @{ Cell = "G2"; Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/b9899ce00336e00ec3646ca43ba553cd838f7f7e/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf"; Display = "075d7a73-6826-41de-8bfe-e6cf8e17f9ec.94eadd793ec9569dbea66e55c42f4ff2303055cc.de-de.xlf" }
        )
    }
)

foreach ($info in $localeSheets) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Drop the stale row.
    $ws.Rows.Item(3).Delete()

    # Refresh the handoff / handback timestamps on the remaining row.
    $ws.Range("E2").Value = $info.HandoffTime
    $ws.Range("H2").Value = $info.HandbackTime

    # Rebuild hyperlinks (Hyperlinks.Delete() clears the whole sheet).
    $ws.Hyperlinks.Delete()
    foreach ($link in $info.Links) {
        $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Address, "", "", $link.Display)
    }
}
